# Rename the "error" terminology to "issue" in the report header row, on
# every worksheet of the workbook (Summary + per-project sheets).
#
# Column E / F of row 10 hold the "Analyzer Errors" / "Scan Errors" header
# strings; replace them with "Analyzer Issues" / "Scan Issues". Since the
# new header text has a (very slightly) different length than the old one,
# Excel's best-fit column E also gets a hair wider, so nudge its width to
# match.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("E10").Value = "Analyzer Issues"
    $ws.Range("F10").Value = "Scan Issues"

    # Column E is sized to best-fit its (now slightly longer) header text.
    $ws.Columns.Item(5).ColumnWidth = 14.17
}
